# Rename the inline picture shapes in the document's headers/footers.
#
# The BTec logo pictures (in the headers) go from "image1.jpg" to "image2.jpg",
# and the Pearson logo pictures (in the footers) go from "image2.png" to
# "image1.png". Every header/footer in every section is scanned so the
# change is applied regardless of how many header/footer slots exist.

$d = $word.ActiveDocument

function Rename-InlineShapes($range) {
    if ($range -eq $null) { return }
    if (-not $range.Exists) { return }
    $shapes = $range.Range.InlineShapes
    if ($shapes.Count -eq 0) { return }
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        $desc = $shp.AlternativeText

        if ($desc -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        } elseif ($desc -like "*PearsonLogo.png") {
            $shp.Name = "image1.png"
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections($s)

    for ($i = 1; $i -le 3; $i++) {
        Rename-InlineShapes($sec.Headers($i))
        Rename-InlineShapes($sec.Footers($i))
    }
}
